$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 18,7
$arr[0,0] = "combination_1_ABCD"
$arr[0,1] = "A"
$arr[0,2] = "KNN"
$arr[0,3] = 10
$arr[0,4] = "AC"
$arr[0,5] = 0.2985953968253966
$arr[0,6] = 1.332453872828025
$arr[1,0] = "combination_1_ABCD"
$arr[1,1] = "A"
$arr[1,2] = "KNN"
$arr[1,3] = 15
$arr[1,4] = "ABCD"
$arr[1,5] = 0.3007955555555555
$arr[1,6] = 1.350484907338554
$arr[2,0] = "combination_1_ABCD"
$arr[2,1] = "A"
$arr[2,2] = "KNN"
$arr[2,3] = 20
$arr[2,4] = "A"
$arr[2,5] = 0.373043180349063
$arr[2,6] = 1.672091178108411
$arr[3,0] = "combination_3_ABCDF"
$arr[3,1] = "B"
$arr[3,2] = "KNN"
$arr[3,3] = 10
$arr[3,4] = "AB"
$arr[3,5] = 0
$arr[3,6] = 0
$arr[4,0] = "combination_3_ABCDF"
$arr[4,1] = "B"
$arr[4,2] = "RandomForest"
$arr[4,3] = 15
$arr[4,4] = "ABC"
$arr[4,5] = 0.0089299999999986
$arr[4,6] = 0.1541971880261283
$arr[5,0] = "combination_3_ABCDF"
$arr[5,1] = "B"
$arr[5,2] = "RandomForest_MICE"
$arr[5,3] = 20
$arr[5,4] = "ABC"
$arr[5,5] = 0.1516413333333328
$arr[5,6] = 2.454008133207574
$arr[6,0] = "combination_3_ABCDF"
$arr[6,1] = "C"
$arr[6,2] = "KNN"
$arr[6,3] = 10
$arr[6,4] = "AC"
$arr[6,5] = 0
$arr[6,6] = 0
$arr[7,0] = "combination_3_ABCDF"
$arr[7,1] = "C"
$arr[7,2] = "RandomForest_MICE"
$arr[7,3] = 15
$arr[7,4] = "ABCD"
$arr[7,5] = 0.0394049999999994
$arr[7,6] = 0.369999999999995
$arr[8,0] = "combination_3_ABCDF"
$arr[8,1] = "C"
$arr[8,2] = "RandomForest_MICE"
$arr[8,3] = 20
$arr[8,4] = "ABCD"
$arr[8,5] = 0.0111479999999996
$arr[8,6] = 0.1046760563380243
$arr[9,0] = "combination_3_ABCDF"
$arr[9,1] = "D"
$arr[9,2] = "KNN"
$arr[9,3] = 10
$arr[9,4] = "AD"
$arr[9,5] = 0
$arr[9,6] = 0
$arr[10,0] = "combination_2_ABCDE"
$arr[10,1] = "D"
$arr[10,2] = "RandomForest"
$arr[10,3] = 15
$arr[10,4] = "BCDE"
$arr[10,5] = 0.0431080714285713
$arr[10,6] = 0.5039291049652861
$arr[11,0] = "combination_3_ABCDF"
$arr[11,1] = "D"
$arr[11,2] = "RandomForest_MICE"
$arr[11,3] = 20
$arr[11,4] = "BCDF"
$arr[11,5] = 0.094764
$arr[11,6] = 1.843657587548639
$arr[12,0] = "combination_2_ABCDE"
$arr[12,1] = "E"
$arr[12,2] = "KNN"
$arr[12,3] = 10
$arr[12,4] = "ACDE"
$arr[12,5] = 0.00024
$arr[12,6] = 0.0014388489208633
$arr[13,0] = "combination_2_ABCDE"
$arr[13,1] = "E"
$arr[13,2] = "KNN"
$arr[13,3] = 15
$arr[13,4] = "BDE"
$arr[13,5] = 0.0223714285714285
$arr[13,6] = 0.1365538104962959
$arr[14,0] = "combination_2_ABCDE"
$arr[14,1] = "E"
$arr[14,2] = "KNN"
$arr[14,3] = 20
$arr[14,4] = "ACE"
$arr[14,5] = 0.0726758457229045
$arr[14,6] = 0.4305782101544178
$arr[15,0] = "combination_3_ABCDF"
$arr[15,1] = "F"
$arr[15,2] = "KNN"
$arr[15,3] = 10
$arr[15,4] = "AF"
$arr[15,5] = 0
$arr[15,6] = 0
$arr[16,0] = "combination_3_ABCDF"
$arr[16,1] = "F"
$arr[16,2] = "HybridKNN_RF"
$arr[16,3] = 15
$arr[16,4] = "BDF"
$arr[16,5] = 0.04698
$arr[16,6] = 0.385081967213115
$arr[17,0] = "combination_3_ABCDF"
$arr[17,1] = "F"
$arr[17,2] = "HybridKNN_RF"
$arr[17,3] = 20
$arr[17,4] = "BDF"
$arr[17,5] = 0.025968
$arr[17,6] = 0.2128524590163936

$rng = $ws.Range("A2:G19")
$rng.Value = $arr
Write-Host "Updated data for Terr Herb Mammals"
